$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 3; $i -le 6; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "Venta $i"
    $ws.Cells.Item($row, 2).Value = $i * 100.0
}
